$d = $word.ActiveDocument

$replacements = @(
    @("51×78=3978", "81×74=5994"),
    @("72×88=6336", "27×11=297"),
    @("88×42=3696", "32×38=1216"),
    @("73×40=2920", "51×94=4794"),
    @("48×60=2880", "75×66=4950"),
    @("11×95=1045", "81×91=7371"),
    @("33×24=792", "29×54=1566"),
    @("59×56=3304", "71×62=4402"),
    @("35×84=2940", "39×26=1014"),
    @("41×16=656", "31×63=1953"),
    @("63×72=4536", "73×28=2044"),
    @("35×97=3395", "48×42=2016"),
    @("88×32=2816", "63×39=2457"),
    @("47×97=4559", "67×65=4355"),
    @("60×34=2040", "14×62=868"),
    @("84×47=3948", "88×80=7040"),
    @("44×16=704", "67×92=6164"),
    @("86×30=2580", "54×51=2754"),
    @("69×12=828", "35×63=2205"),
    @("53×18=954", "64×24=1536"),
    @("81×96=7776", "23×81=1863"),
    @("25×76=1900", "34×31=1054"),
    @("65×58=3770", "27×57=1539"),
    @("85×76=6460", "18×39=702"),
    @("40×82=3280", "57×17=969")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
